$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "3"
$ws.Range("G2").Value = "109.1118546666667"
$ws.Range("H2").Value = "327.335564"
$ws.Range("I2").Value = "0.3029068882986101"
$ws.Range("J2").Value = "0.3029068882986101"
$ws.Range("K2").Value = "3"
$ws.Range("M2").Value = "2.092292333333333"
$ws.Range("N2").Value = "6.276877000000001"
$ws.Range("O2").Value = "0.1237967521619938"
$ws.Range("P2").Value = "0.1237967521619938"
$ws.Range("Q2").Value = "228.2938969948476"
$ws.Range("R2").Value = "2054.645072953629"
$ws.Range("S2").Value = "0.03749888897886376"
$ws.Range("T2").Value = "0.03749888897886377"

$ws.Range("E3").Value = "3"
$ws.Range("G3").Value = "109.1118546666667"
$ws.Range("H3").Value = "327.335564"
$ws.Range("I3").Value = "0.3029068882986101"
$ws.Range("J3").Value = "0.3029068882986101"
$ws.Range("K3").Value = "3"
$ws.Range("M3").Value = "2.468365333333333"
$ws.Range("N3").Value = "7.405096"
$ws.Range("O3").Value = "0.1460482393151517"
$ws.Range("P3").Value = "0.1460482393151517"
$ws.Range("Q3").Value = "269.3279195149049"
$ws.Range("R3").Value = "2423.951275634144"
$ws.Range("S3").Value = "0.04423901771244333"
$ws.Range("T3").Value = "0.04423901771244333"

$ws.Range("E4").Value = "3"
$ws.Range("G4").Value = "109.1118546666667"
$ws.Range("H4").Value = "327.335564"
$ws.Range("I4").Value = "0.3029068882986101"
$ws.Range("J4").Value = "0.3029068882986101"
$ws.Range("K4").Value = "3"
$ws.Range("M4").Value = "0.315935"
$ws.Range("N4").Value = "0.9478049999999999"
$ws.Range("O4").Value = "0.01869324198688273"
$ws.Range("P4").Value = "0.01869324198688273"
$ws.Range("Q4").Value = "34.47225380411333"
$ws.Range("R4").Value = "310.25028423702"
$ws.Range("S4").Value = "0.005662311762459574"
$ws.Range("T4").Value = "0.005662311762459574"

$ws.Range("E5").Value = "3"
$ws.Range("G5").Value = "109.1118546666667"
$ws.Range("H5").Value = "327.335564"
$ws.Range("I5").Value = "0.3029068882986101"
$ws.Range("J5").Value = "0.3029068882986101"
$ws.Range("K5").Value = "3"
$ws.Range("M5").Value = "0.4705663333333334"
$ws.Range("N5").Value = "1.411699"
$ws.Range("O5").Value = "0.02784246867197405"
$ws.Range("P5").Value = "0.02784246867197405"
$ws.Range("Q5").Value = "51.3443653736929"
$ws.Range("R5").Value = "462.099288363236"
$ws.Range("S5").Value = "0.008433675547979193"
$ws.Range("T5").Value = "0.008433675547979195"

$ws.Range("E6").Value = "3"
$ws.Range("G6").Value = "109.1118546666667"
$ws.Range("H6").Value = "327.335564"
$ws.Range("I6").Value = "0.3029068882986101"
$ws.Range("J6").Value = "0.3029068882986101"
$ws.Range("K6").Value = "3"
$ws.Range("M6").Value = "11.291786"
$ws.Range("N6").Value = "33.875358"
$ws.Range("O6").Value = "0.6681123907199095"
$ws.Range("P6").Value = "0.6681123907199095"
$ws.Range("Q6").Value = "1232.067712959102"
$ws.Range("R6").Value = "11088.60941663191"
$ws.Range("S6").Value = "0.2023758453067129"
$ws.Range("T6").Value = "0.2023758453067129"

$ws.Range("E7").Value = "3"
$ws.Range("G7").Value = "109.1118546666667"
$ws.Range("H7").Value = "327.335564"
$ws.Range("I7").Value = "0.3029068882986101"
$ws.Range("J7").Value = "0.3029068882986101"
$ws.Range("K7").Value = "3"
$ws.Range("M7").Value = "0.2620826666666667"
$ws.Range("N7").Value = "0.7862480000000001"
$ws.Range("O7").Value = "0.01550690714408826"
$ws.Range("P7").Value = "0.01550690714408826"
$ws.Range("Q7").Value = "28.59632583598578"
$ws.Range("R7").Value = "257.366932523872"
$ws.Range("S7").Value = "0.004697148990151261"
$ws.Range("T7").Value = "0.004697148990151261"

$ws.Range("E8").Value = "3"
$ws.Range("G8").Value = "227.6338753333333"
$ws.Range("H8").Value = "682.901626"
$ws.Range("I8").Value = "0.6319374650831437"
$ws.Range("J8").Value = "0.6319374650831437"
$ws.Range("K8").Value = "3"
$ws.Range("M8").Value = "2.092292333333333"
$ws.Range("N8").Value = "6.276877000000001"
$ws.Range("O8").Value = "0.1237967521619938"
$ws.Range("P8").Value = "0.1237967521619938"
$ws.Range("Q8").Value = "476.2766121668891"
$ws.Range("R8").Value = "4286.489509502002"
$ws.Range("S8").Value = "0.07823180574677654"
$ws.Range("T8").Value = "0.07823180574677654"

$ws.Range("E9").Value = "3"
$ws.Range("G9").Value = "227.6338753333333"
$ws.Range("H9").Value = "682.901626"
$ws.Range("I9").Value = "0.6319374650831437"
$ws.Range("J9").Value = "0.6319374650831437"
$ws.Range("K9").Value = "3"
$ws.Range("M9").Value = "2.468365333333333"
$ws.Range("N9").Value = "7.405096"
$ws.Range("O9").Value = "0.1460482393151517"
$ws.Range("P9").Value = "0.1460482393151517"
$ws.Range("Q9").Value = "561.8835665651218"
$ws.Range("R9").Value = "5056.952099086096"
$ws.Range("S9").Value = "0.0922933541326733"
$ws.Range("T9").Value = "0.0922933541326733"

$ws.Range("E10").Value = "3"
$ws.Range("G10").Value = "227.6338753333333"
$ws.Range("H10").Value = "682.901626"
$ws.Range("I10").Value = "0.6319374650831437"
$ws.Range("J10").Value = "0.6319374650831437"
$ws.Range("K10").Value = "3"
$ws.Range("M10").Value = "0.315935"
$ws.Range("N10").Value = "0.9478049999999999"
$ws.Range("O10").Value = "0.01869324198688273"
$ws.Range("P10").Value = "0.01869324198688273"
$ws.Range("Q10").Value = "71.91750840343666"
$ws.Range("R10").Value = "647.2575756309299"
$ws.Range("S10").Value = "0.01181295995537646"
$ws.Range("T10").Value = "0.01181295995537646"

$ws.Range("E11").Value = "3"
$ws.Range("G11").Value = "227.6338753333333"
$ws.Range("H11").Value = "682.901626"
$ws.Range("I11").Value = "0.6319374650831437"
$ws.Range("J11").Value = "0.6319374650831437"
$ws.Range("K11").Value = "3"
$ws.Range("M11").Value = "0.4705663333333334"
$ws.Range("N11").Value = "1.411699"
$ws.Range("O11").Value = "0.02784246867197405"
$ws.Range("P11").Value = "0.02784246867197405"
$ws.Range("Q11").Value = "107.1168380580638"
$ws.Range("R11").Value = "964.051542522574"
$ws.Range("S11").Value = "0.01759469907422413"
$ws.Range("T11").Value = "0.01759469907422413"

$ws.Range("E12").Value = "3"
$ws.Range("G12").Value = "227.6338753333333"
$ws.Range("H12").Value = "682.901626"
$ws.Range("I12").Value = "0.6319374650831437"
$ws.Range("J12").Value = "0.6319374650831437"
$ws.Range("K12").Value = "3"
$ws.Range("M12").Value = "11.291786"
$ws.Range("N12").Value = "33.875358"
$ws.Range("O12").Value = "0.6681123907199095"
$ws.Range("P12").Value = "0.6681123907199095"
$ws.Range("Q12").Value = "2570.393006614679"
$ws.Range("R12").Value = "23133.53705953211"
$ws.Range("S12").Value = "0.4222052505821785"
$ws.Range("T12").Value = "0.4222052505821785"

$ws.Range("E13").Value = "3"
$ws.Range("G13").Value = "227.6338753333333"
$ws.Range("H13").Value = "682.901626"
$ws.Range("I13").Value = "0.6319374650831437"
$ws.Range("J13").Value = "0.6319374650831437"
$ws.Range("K13").Value = "3"
$ws.Range("M13").Value = "0.2620826666666667"
$ws.Range("N13").Value = "0.7862480000000001"
$ws.Range("O13").Value = "0.01550690714408826"
$ws.Range("P13").Value = "0.01550690714408826"
$ws.Range("Q13").Value = "59.65889307102756"
$ws.Range("R13").Value = "536.930037639248"
$ws.Range("S13").Value = "0.009799395591914827"
$ws.Range("T13").Value = "0.009799395591914827"

$ws.Range("E14").Value = "3"
$ws.Range("G14").Value = "0.194568"
$ws.Range("H14").Value = "0.583704"
$ws.Range("I14").Value = "0.0005401428435299865"
$ws.Range("J14").Value = "0.0005401428435299865"
$ws.Range("K14").Value = "3"
$ws.Range("M14").Value = "2.092292333333333"
$ws.Range("N14").Value = "6.276877000000001"
$ws.Range("O14").Value = "0.1237967521619938"
$ws.Range("P14").Value = "0.1237967521619938"
$ws.Range("Q14").Value = "0.407093134712"
$ws.Range("R14").Value = "3.663838212408"
$ws.Range("S14").Value = "6.686792973255632E-05"
$ws.Range("T14").Value = "6.686792973255633E-05"

$ws.Range("E15").Value = "3"
$ws.Range("G15").Value = "0.194568"
$ws.Range("H15").Value = "0.583704"
$ws.Range("I15").Value = "0.0005401428435299865"
$ws.Range("J15").Value = "0.0005401428435299865"
$ws.Range("K15").Value = "3"
$ws.Range("M15").Value = "2.468365333333333"
$ws.Range("N15").Value = "7.405096"
$ws.Range("O15").Value = "0.1460482393151517"
$ws.Range("P15").Value = "0.1460482393151517"
$ws.Range("Q15").Value = "0.480264906176"
$ws.Range("R15").Value = "4.322384155584"
$ws.Range("S15").Value = "7.888691127623401E-05"
$ws.Range("T15").Value = "7.888691127623401E-05"

$ws.Range("E16").Value = "3"
$ws.Range("G16").Value = "0.194568"
$ws.Range("H16").Value = "0.583704"
$ws.Range("I16").Value = "0.0005401428435299865"
$ws.Range("J16").Value = "0.0005401428435299865"
$ws.Range("K16").Value = "3"
$ws.Range("M16").Value = "0.315935"
$ws.Range("N16").Value = "0.9478049999999999"
$ws.Range("O16").Value = "0.01869324198688273"
$ws.Range("P16").Value = "0.01869324198688273"
$ws.Range("Q16").Value = "0.06147084107999999"
$ws.Range("R16").Value = "0.5532375697199999"
$ws.Range("S16").Value = "1.009702088158897E-05"
$ws.Range("T16").Value = "1.009702088158897E-05"

$ws.Range("E17").Value = "3"
$ws.Range("G17").Value = "0.194568"
$ws.Range("H17").Value = "0.583704"
$ws.Range("I17").Value = "0.0005401428435299865"
$ws.Range("J17").Value = "0.0005401428435299865"
$ws.Range("K17").Value = "3"
$ws.Range("M17").Value = "0.4705663333333334"
$ws.Range("N17").Value = "1.411699"
$ws.Range("O17").Value = "0.02784246867197405"
$ws.Range("P17").Value = "0.02784246867197405"
$ws.Range("Q17").Value = "0.091557150344"
$ws.Range("R17").Value = "0.824014353096"
$ws.Range("S17").Value = "1.503891019937463E-05"
$ws.Range("T17").Value = "1.503891019937463E-05"

$ws.Range("E18").Value = "3"
$ws.Range("G18").Value = "0.194568"
$ws.Range("H18").Value = "0.583704"
$ws.Range("I18").Value = "0.0005401428435299865"
$ws.Range("J18").Value = "0.0005401428435299865"
$ws.Range("K18").Value = "3"
$ws.Range("M18").Value = "11.291786"
$ws.Range("N18").Value = "33.875358"
$ws.Range("O18").Value = "0.6681123907199095"
$ws.Range("P18").Value = "0.6681123907199095"
$ws.Range("Q18").Value = "2.197020218448"
$ws.Range("R18").Value = "19.773181966032"
$ws.Range("S18").Value = "0.0003608761265210693"
$ws.Range("T18").Value = "0.0003608761265210693"

$ws.Range("E19").Value = "3"
$ws.Range("G19").Value = "0.194568"
$ws.Range("H19").Value = "0.583704"
$ws.Range("I19").Value = "0.0005401428435299865"
$ws.Range("J19").Value = "0.0005401428435299865"
$ws.Range("K19").Value = "3"
$ws.Range("M19").Value = "0.2620826666666667"
$ws.Range("N19").Value = "0.7862480000000001"
$ws.Range("O19").Value = "0.01550690714408826"
$ws.Range("P19").Value = "0.01550690714408826"
$ws.Range("Q19").Value = "0.050992900288"
$ws.Range("R19").Value = "0.458936102592"
$ws.Range("S19").Value = "8.375944919163296E-06"
$ws.Range("T19").Value = "8.375944919163296E-06"

$ws.Range("E20").Value = "3"
$ws.Range("G20").Value = "0.4517016666666667"
$ws.Range("H20").Value = "1.355105"
$ws.Range("I20").Value = "0.001253975076377243"
$ws.Range("J20").Value = "0.001253975076377243"
$ws.Range("K20").Value = "3"
$ws.Range("M20").Value = "2.092292333333333"
$ws.Range("N20").Value = "6.276877000000001"
$ws.Range("O20").Value = "0.1237967521619938"
$ws.Range("P20").Value = "0.1237967521619938"
$ws.Range("Q20").Value = "0.9450919341205556"
$ws.Range("R20").Value = "8.505827407085"
$ws.Range("S20").Value = "0.0001552380417475908"
$ws.Range("T20").Value = "0.0001552380417475908"

$ws.Range("E21").Value = "3"
$ws.Range("G21").Value = "0.4517016666666667"
$ws.Range("H21").Value = "1.355105"
$ws.Range("I21").Value = "0.001253975076377243"
$ws.Range("J21").Value = "0.001253975076377243"
$ws.Range("K21").Value = "3"
$ws.Range("M21").Value = "2.468365333333333"
$ws.Range("N21").Value = "7.405096"
$ws.Range("O21").Value = "0.1460482393151517"
$ws.Range("P21").Value = "0.1460482393151517"
$ws.Range("Q21").Value = "1.114964735008889"
$ws.Range("R21").Value = "10.03468261508"
$ws.Range("S21").Value = "0.0001831408520499792"
$ws.Range("T21").Value = "0.0001831408520499792"

$ws.Range("E22").Value = "3"
$ws.Range("G22").Value = "0.4517016666666667"
$ws.Range("H22").Value = "1.355105"
$ws.Range("I22").Value = "0.001253975076377243"
$ws.Range("J22").Value = "0.001253975076377243"
$ws.Range("K22").Value = "3"
$ws.Range("M22").Value = "0.315935"
$ws.Range("N22").Value = "0.9478049999999999"
$ws.Range("O22").Value = "0.01869324198688273"
$ws.Range("P22").Value = "0.01869324198688273"
$ws.Range("Q22").Value = "0.1427083660583333"
$ws.Range("R22").Value = "1.284375294525"
$ws.Range("S22").Value = "2.344085954823956E-05"
$ws.Range("T22").Value = "2.344085954823956E-05"

$ws.Range("E23").Value = "3"
$ws.Range("G23").Value = "0.4517016666666667"
$ws.Range("H23").Value = "1.355105"
$ws.Range("I23").Value = "0.001253975076377243"
$ws.Range("J23").Value = "0.001253975076377243"
$ws.Range("K23").Value = "3"
$ws.Range("M23").Value = "0.4705663333333334"
$ws.Range("N23").Value = "1.411699"
$ws.Range("O23").Value = "0.02784246867197405"
$ws.Range("P23").Value = "0.02784246867197405"
$ws.Range("Q23").Value = "0.2125555970438889"
$ws.Range("R23").Value = "1.913000373395"
$ws.Range("S23").Value = "3.491376177946965E-05"
$ws.Range("T23").Value = "3.491376177946966E-05"

$ws.Range("E24").Value = "3"
$ws.Range("G24").Value = "0.4517016666666667"
$ws.Range("H24").Value = "1.355105"
$ws.Range("I24").Value = "0.001253975076377243"
$ws.Range("J24").Value = "0.001253975076377243"
$ws.Range("K24").Value = "3"
$ws.Range("M24").Value = "11.291786"
$ws.Range("N24").Value = "33.875358"
$ws.Range("O24").Value = "0.6681123907199095"
$ws.Range("P24").Value = "0.6681123907199095"
$ws.Range("Q24").Value = "5.100518555843333"
$ws.Range("R24").Value = "45.90466700259"
$ws.Range("S24").Value = "0.000837796286181581"
$ws.Range("T24").Value = "0.000837796286181581"

$ws.Range("E25").Value = "3"
$ws.Range("G25").Value = "0.4517016666666667"
$ws.Range("H25").Value = "1.355105"
$ws.Range("I25").Value = "0.001253975076377243"
$ws.Range("J25").Value = "0.001253975076377243"
$ws.Range("K25").Value = "3"
$ws.Range("M25").Value = "0.2620826666666667"
$ws.Range("N25").Value = "0.7862480000000001"
$ws.Range("O25").Value = "0.01550690714408826"
$ws.Range("P25").Value = "0.01550690714408826"
$ws.Range("Q25").Value = "0.1183831773377778"
$ws.Range("R25").Value = "1.06544859604"
$ws.Range("S25").Value = "1.944527507038289E-05"
$ws.Range("T25").Value = "1.944527507038289E-05"

$ws.Range("E26").Value = "3"
$ws.Range("G26").Value = "0.7729569999999999"
$ws.Range("H26").Value = "2.318871"
$ws.Range("I26").Value = "0.002145816331084288"
$ws.Range("J26").Value = "0.002145816331084288"
$ws.Range("K26").Value = "3"
$ws.Range("M26").Value = "2.092292333333333"
$ws.Range("N26").Value = "6.276877000000001"
$ws.Range("O26").Value = "0.1237967521619938"
$ws.Range("P26").Value = "0.1237967521619938"
$ws.Range("Q26").Value = "1.617252005096333"
$ws.Range("R26").Value = "14.555268045867"
$ws.Range("S26").Value = "0.0002656450925244004"
$ws.Range("T26").Value = "0.0002656450925244004"

$ws.Range("E27").Value = "3"
$ws.Range("G27").Value = "0.7729569999999999"
$ws.Range("H27").Value = "2.318871"
$ws.Range("I27").Value = "0.002145816331084288"
$ws.Range("J27").Value = "0.002145816331084288"
$ws.Range("K27").Value = "3"
$ws.Range("M27").Value = "2.468365333333333"
$ws.Range("N27").Value = "7.405096"
$ws.Range("O27").Value = "0.1460482393151517"
$ws.Range("P27").Value = "0.1460482393151517"
$ws.Range("Q27").Value = "1.907940262957333"
$ws.Range("R27").Value = "17.171462366616"
$ws.Range("S27").Value = "0.0003133926970485589"
$ws.Range("T27").Value = "0.0003133926970485589"

$ws.Range("E28").Value = "3"
$ws.Range("G28").Value = "0.7729569999999999"
$ws.Range("H28").Value = "2.318871"
$ws.Range("I28").Value = "0.002145816331084288"
$ws.Range("J28").Value = "0.002145816331084288"
$ws.Range("K28").Value = "3"
$ws.Range("M28").Value = "0.315935"
$ws.Range("N28").Value = "0.9478049999999999"
$ws.Range("O28").Value = "0.01869324198688273"
$ws.Range("P28").Value = "0.01869324198688273"
$ws.Range("Q28").Value = "0.244204169795"
$ws.Range("R28").Value = "2.197837528155"
$ws.Range("S28").Value = "4.011226393636347E-05"
$ws.Range("T28").Value = "4.011226393636347E-05"

$ws.Range("E29").Value = "3"
$ws.Range("G29").Value = "0.7729569999999999"
$ws.Range("H29").Value = "2.318871"
$ws.Range("I29").Value = "0.002145816331084288"
$ws.Range("J29").Value = "0.002145816331084288"
$ws.Range("K29").Value = "3"
$ws.Range("M29").Value = "0.4705663333333334"
$ws.Range("N29").Value = "1.411699"
$ws.Range("O29").Value = "0.02784246867197405"
$ws.Range("P29").Value = "0.02784246867197405"
$ws.Range("Q29").Value = "0.3637275413143333"
$ws.Range("R29").Value = "3.273547871828999"
$ws.Range("S29").Value = "5.974482397402458E-05"
$ws.Range("T29").Value = "5.974482397402459E-05"

$ws.Range("E30").Value = "3"
$ws.Range("G30").Value = "0.7729569999999999"
$ws.Range("H30").Value = "2.318871"
$ws.Range("I30").Value = "0.002145816331084288"
$ws.Range("J30").Value = "0.002145816331084288"
$ws.Range("K30").Value = "3"
$ws.Range("M30").Value = "11.291786"
$ws.Range("N30").Value = "33.875358"
$ws.Range("O30").Value = "0.6681123907199095"
$ws.Range("P30").Value = "0.6681123907199095"
$ws.Range("Q30").Value = "8.728065031201998"
$ws.Range("R30").Value = "78.55258528081798"
$ws.Range("S30").Value = "0.001433646479006548"
$ws.Range("T30").Value = "0.001433646479006548"

$ws.Range("E31").Value = "3"
$ws.Range("G31").Value = "0.7729569999999999"
$ws.Range("H31").Value = "2.318871"
$ws.Range("I31").Value = "0.002145816331084288"
$ws.Range("J31").Value = "0.002145816331084288"
$ws.Range("K31").Value = "3"
$ws.Range("M31").Value = "0.2620826666666667"
$ws.Range("N31").Value = "0.7862480000000001"
$ws.Range("O31").Value = "0.01550690714408826"
$ws.Range("P31").Value = "0.01550690714408826"
$ws.Range("Q31").Value = "0.2025786317786666"
$ws.Range("R31").Value = "1.823207686008"
$ws.Range("S31").Value = "3.327497459439221E-05"
$ws.Range("T31").Value = "3.327497459439221E-05"

$ws.Range("E32").Value = "3"
$ws.Range("G32").Value = "22.05086833333333"
$ws.Range("H32").Value = "66.15260499999999"
$ws.Range("I32").Value = "0.06121571236725463"
$ws.Range("J32").Value = "0.06121571236725463"
$ws.Range("K32").Value = "3"
$ws.Range("M32").Value = "2.092292333333333"
$ws.Range("N32").Value = "6.276877000000001"
$ws.Range("O32").Value = "0.1237967521619938"
$ws.Range("P32").Value = "0.1237967521619938"
$ws.Range("Q32").Value = "46.1368627571761"
$ws.Range("R32").Value = "415.231764814585"
$ws.Range("S32").Value = "0.007578306372348919"
$ws.Range("T32").Value = "0.00757830637234892"

$ws.Range("E33").Value = "3"
$ws.Range("G33").Value = "22.05086833333333"
$ws.Range("H33").Value = "66.15260499999999"
$ws.Range("I33").Value = "0.06121571236725463"
$ws.Range("J33").Value = "0.06121571236725463"
$ws.Range("K33").Value = "3"
$ws.Range("M33").Value = "2.468365333333333"
$ws.Range("N33").Value = "7.405096"
$ws.Range("O33").Value = "0.1460482393151517"
$ws.Range("P33").Value = "0.1460482393151517"
$ws.Range("Q33").Value = "54.42959896389777"
$ws.Range("R33").Value = "489.86639067508"
$ws.Range("S33").Value = "0.008940447009660298"
$ws.Range("T33").Value = "0.008940447009660298"

$ws.Range("E34").Value = "3"
$ws.Range("G34").Value = "22.05086833333333"
$ws.Range("H34").Value = "66.15260499999999"
$ws.Range("I34").Value = "0.06121571236725463"
$ws.Range("J34").Value = "0.06121571236725463"
$ws.Range("K34").Value = "3"
$ws.Range("M34").Value = "0.315935"
$ws.Range("N34").Value = "0.9478049999999999"
$ws.Range("O34").Value = "0.01869324198688273"
$ws.Range("P34").Value = "0.01869324198688273"
$ws.Range("Q34").Value = "6.966641086891665"
$ws.Range("R34").Value = "62.69976978202499"
$ws.Range("S34").Value = "0.001144320124680501"
$ws.Range("T34").Value = "0.001144320124680501"

$ws.Range("E35").Value = "3"
$ws.Range("G35").Value = "22.05086833333333"
$ws.Range("H35").Value = "66.15260499999999"
$ws.Range("I35").Value = "0.06121571236725463"
$ws.Range("J35").Value = "0.06121571236725463"
$ws.Range("K35").Value = "3"
$ws.Range("M35").Value = "0.4705663333333334"
$ws.Range("N35").Value = "1.411699"
$ws.Range("O35").Value = "0.02784246867197405"
$ws.Range("P35").Value = "0.02784246867197405"
$ws.Range("Q35").Value = "10.37639625843278"
$ws.Range("R35").Value = "93.38756632589499"
$ws.Range("S35").Value = "0.001704396553817862"
$ws.Range("T35").Value = "0.001704396553817862"

$ws.Range("E36").Value = "3"
$ws.Range("G36").Value = "22.05086833333333"
$ws.Range("H36").Value = "66.15260499999999"
$ws.Range("I36").Value = "0.06121571236725463"
$ws.Range("J36").Value = "0.06121571236725463"
$ws.Range("K36").Value = "3"
$ws.Range("M36").Value = "11.291786"
$ws.Range("N36").Value = "33.875358"
$ws.Range("O36").Value = "0.6681123907199095"
$ws.Range("P36").Value = "0.6681123907199095"
$ws.Range("Q36").Value = "248.9936863341766"
$ws.Range("R36").Value = "2240.94317700759"
$ws.Range("S36").Value = "0.04089897593930882"
$ws.Range("T36").Value = "0.04089897593930882"

$ws.Range("E37").Value = "3"
$ws.Range("G37").Value = "22.05086833333333"
$ws.Range("H37").Value = "66.15260499999999"
$ws.Range("I37").Value = "0.06121571236725463"
$ws.Range("J37").Value = "0.06121571236725463"
$ws.Range("K37").Value = "3"
$ws.Range("M37").Value = "0.2620826666666667"
$ws.Range("N37").Value = "0.7862480000000001"
$ws.Range("O37").Value = "0.01550690714408826"
$ws.Range("P37").Value = "0.01550690714408826"
$ws.Range("Q37").Value = "5.779150375115555"
$ws.Range("R37").Value = "52.01235337604"
$ws.Range("S37").Value = "0.000949266367438233"
$ws.Range("T37").Value = "0.000949266367438233"
